$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name, Link) - straightforward text, no numeric ambiguity
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

# Price / Volume columns - force text storage to avoid Excel auto-numeric conversion,
# then reset style back to Normal so no extra style index is left on the cell.
$priceVolumeCells = @{}
$priceVolumeCells['D2'] = '51.630.58'
$priceVolumeCells['E2'] = '  -0.88%  '
$priceVolumeCells['D3'] = '2.785.05'
$priceVolumeCells['E3'] = '  -0.38%  '
$priceVolumeCells['E4'] = '  -0.01%  '
$priceVolumeCells['D5'] = '352.43'
$priceVolumeCells['E5'] = '  -2.14%  '
$priceVolumeCells['D6'] = '109.09'
$priceVolumeCells['E6'] = '  -0.97%  '
$priceVolumeCells['D7'] = '0.551'
$priceVolumeCells['E7'] = '  -2.35%  '
$priceVolumeCells['E8'] = '  +0.03%  '
$priceVolumeCells['D9'] = '0.609'
$priceVolumeCells['E9'] = '  +2.08%  '
$priceVolumeCells['E10'] = '  -1.42%  '
$priceVolumeCells['E11'] = '  +2.53%  '
$priceVolumeCells['D13'] = '20.09'
$priceVolumeCells['E13'] = '  +2.83%  '
$priceVolumeCells['D14'] = '7.69'
$priceVolumeCells['E14'] = '  +0.82%  '
$priceVolumeCells['D15'] = '3.219.88'
$priceVolumeCells['E15'] = '  -0.42%  '
$priceVolumeCells['D16'] = '2.786.48'
$priceVolumeCells['E16'] = '  -0.35%  '
$priceVolumeCells['E17'] = '  -2.40%  '
$priceVolumeCells['D18'] = '51.603.46'
$priceVolumeCells['E18'] = '  -0.73%  '
$priceVolumeCells['D19'] = '7.69'
$priceVolumeCells['E19'] = '  +3.55%  '
$priceVolumeCells['E20'] = '  +0.54%  '
$priceVolumeCells['D21'] = '13.17'
$priceVolumeCells['E21'] = '  +0.77%  '
$priceVolumeCells['D22'] = '0.0₃0964'
$priceVolumeCells['E22'] = '  -2.33%  '
$priceVolumeCells['D23'] = '69.91'
$priceVolumeCells['E23'] = '  -0.64%  '
$priceVolumeCells['D24'] = '267.15'
$priceVolumeCells['E24'] = '  -2.59%  '
$priceVolumeCells['D25'] = '2.73'
$priceVolumeCells['E25'] = '  -1.01%  '
$priceVolumeCells['E26'] = '  -2.30%  '
$priceVolumeCells['E27'] = '  -0.18%  '
$priceVolumeCells['D28'] = '0.165'
$priceVolumeCells['E28'] = '  +12.85%  '
$priceVolumeCells['D29'] = '10.25'
$priceVolumeCells['E29'] = '  +0.39%  '
$priceVolumeCells['D30'] = '36.93'
$priceVolumeCells['E30'] = '  +7.08%  '
$priceVolumeCells['E31'] = '  -2.34%  '
$priceVolumeCells['E32'] = '  +6.88%  '
$priceVolumeCells['D33'] = '51.70'
$priceVolumeCells['E33'] = '  +0.03%  '
$priceVolumeCells['D34'] = '5.71'
$priceVolumeCells['E34'] = '  +8.50%  '
$priceVolumeCells['D35'] = '0.0454'
$priceVolumeCells['E35'] = '  -2.06%  '
$priceVolumeCells['D36'] = '0.0834'
$priceVolumeCells['E36'] = '  -2.00%  '
$priceVolumeCells['E37'] = '  +0.01%  '
$priceVolumeCells['D38'] = '18.54'
$priceVolumeCells['E38'] = '  +0.98%  '
$priceVolumeCells['E39'] = '  -3.01%  '
$priceVolumeCells['E40'] = '  -1.91%  '
$priceVolumeCells['E41'] = '  -1.51%  '
$priceVolumeCells['E42'] = '  -0.80%  '
$priceVolumeCells['D43'] = '120.29'
$priceVolumeCells['E43'] = '  -2.12%  '
$priceVolumeCells['D44'] = '22.09'
$priceVolumeCells['E44'] = '  -0.79%  '
$priceVolumeCells['D45'] = '2.18'
$priceVolumeCells['E45'] = '  -3.30%  '
$priceVolumeCells['D46'] = '2.126.07'
$priceVolumeCells['E46'] = '  +1.96%  '
$priceVolumeCells['D47'] = '3.33'
$priceVolumeCells['E47'] = '  +1.64%  '
$priceVolumeCells['E48'] = '  +4.27%  '
$priceVolumeCells['D49'] = '5.44'
$priceVolumeCells['E49'] = '  -5.10%  '
$priceVolumeCells['D50'] = '0.906'
$priceVolumeCells['E50'] = '  -3.36%  '
$priceVolumeCells['E51'] = '  +8.03%  '

foreach ($addr in $priceVolumeCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceVolumeCells[$addr]
    $cell.Style = "Normal"
}
